# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" headers on the two existing sheets.
# 2. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header style (bold, bordered, centered) from the
# "Weekly Quantity" sheet's header row, and the date style from its date
# column, so no duplicate styles get created.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$dates = @(45501.99999999999, 45515.99999999999, 45522.99999999999, 45529.99999999999, 45536.99999999999, 45543.99999999999, 45550.99999999999, 45557.99999999999, 45564.99999999999, 45571.99999999999, 45578.99999999999)
$forecast = @(16, 16, 16, 16, 16, 16, 16, 16, 16, 16, 16)
$lower = @(15.99999998068972, 15.99999998093721, 15.99999997937375, 15.99999997372575, 15.99999993986004, 15.99999987604749, 15.99999980487363, 15.99999971278722, 15.99999962010924, 15.99999950064243, 15.9999993627494)
$upper = @(16.00000002018895, 16.00000002079527, 16.00000002141319, 16.00000002650836, 16.00000005982204, 16.00000012107799, 16.00000020921176, 16.00000029792898, 16.00000040745257, 16.00000053069045, 16.00000064835141)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value = $forecast[$i]
    $wsForecast.Cells.Item($row, 3).Value = $lower[$i]
    $wsForecast.Cells.Item($row, 4).Value = $upper[$i]
}
